$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Insert two new rows above the first "Notes:" bullet (old row 10), shifting
# the existing notes (and everything below them) down by two rows.
$ws.Range("A10:A11").EntireRow.Insert()

# Row 10 becomes a brand-new note explaining that the EU EPS reuses US data.
$ws.Range("A10").Value = "The EU EPS uses values from the US EPS."
# Match the formatting of the other note lines (normal, non-bold text) --
# the newly inserted row otherwise inherits the bold "Notes:" header style.
$ws.Range("A10").Font.Bold = $false

# Row 11 is left blank, matching the blank spacer rows used elsewhere on
# this sheet (e.g. between the notes list and the next section).

Write-Host "Inserted new note row in About sheet"
